# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos.xlsx price/volume update described by the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-looking Price (column D) values must be forced to Text format
# --- so Excel stores the exact original string instead of converting it to
# --- a floating point number (which would introduce binary rounding noise
# --- and strip significant trailing/leading zeros).
$numericPriceCells = @(
    "D4", "D5", "D7", "D8", "D10", "D13", "D14", "D18", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D32", "D37", "D39", "D40", "D42", "D43", "D44", "D46", "D47", "D50"
)
foreach ($cellRef in $numericPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Apply the new values
$ws.Range("D2").Value = '30.811.27'
$ws.Range("D3").Value = '1.691.88'
$ws.Range("E3").Value = '  +2.95%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '221.79'
$ws.Range("E5").Value = '  +2.56%  '
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '30.94'
$ws.Range("E8").Value = '  +4.95%  '
$ws.Range("E9").Value = '  +1.71%  '
$ws.Range("D10").Value = '0.0627'
$ws.Range("E10").Value = '  +1.90%  '
$ws.Range("E11").Value = '  -1.36%  '
$ws.Range("D12").Value = '1.937.02'
$ws.Range("E12").Value = '  +3.15%  '
$ws.Range("D13").Value = '10.66'
$ws.Range("E13").Value = '  +11.14%  '
$ws.Range("D14").Value = '0.626'
$ws.Range("E14").Value = '  +8.13%  '
$ws.Range("D15").Value = '1.701.62'
$ws.Range("E15").Value = '  +3.58%  '
$ws.Range("E16").Value = '  +2.62%  '
$ws.Range("D17").Value = '30.844.81'
$ws.Range("E17").Value = '  +2.30%  '
$ws.Range("D18").Value = '66.50'
$ws.Range("E18").Value = '  +2.24%  '
$ws.Range("D19").Value = '248.35'
$ws.Range("E19").Value = '  -0.16%  '
$ws.Range("E20").Value = '  +1.34%  '
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = '10.36'
$ws.Range("E22").Value = '  +3.55%  '
$ws.Range("E23").Value = '  +1.55%  '
$ws.Range("D24").Value = '2.16'
$ws.Range("E24").Value = '  +0.89%  '
$ws.Range("D25").Value = '157.42'
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("D26").Value = '15.93'
$ws.Range("E26").Value = '  +0.87%  '
$ws.Range("D27").Value = '0.112'
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("E28").Value = '  +0.69%  '
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +1.48%  '
$ws.Range("E31").Value = '  +1.15%  '
$ws.Range("D32").Value = '3.49'
$ws.Range("E32").Value = '  +2.04%  '
$ws.Range("D33").Value = '1.520.86'
$ws.Range("E33").Value = '  +5.78%  '
$ws.Range("E34").Value = '  +2.83%  '
$ws.Range("E35").Value = '  +4.35%  '
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").Value = '84.73'
$ws.Range("E37").Value = '  +8.71%  '
$ws.Range("E38").Value = '  +4.47%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.73'
$ws.Range("E39").Value = '  -4.95%  '
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").Value = '0.586'
$ws.Range("E40").Value = '  +4.09%  '
$ws.Range("E41").Value = '  +1.26%  '
$ws.Range("D42").Value = '0.855'
$ws.Range("E42").Value = '  +1.28%  '
$ws.Range("D43").Value = '2.02'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").Value = '0.0504'
$ws.Range("E44").Value = '  +0.98%  '
$ws.Range("E45").Value = '  -1.51%  '
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").Value = '52.54'
$ws.Range("E47").Value = '  -5.80%  '
$ws.Range("D48").Value = '1.830.27'
$ws.Range("E48").Value = '  +2.50%  '
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("D50").Value = '95.49'
$ws.Range("E50").Value = '  +5.51%  '
$ws.Range("D51").Value = '0.0₆0115'
$ws.Range("E51").Value = '  +4.31%  '
